# Updated Dam Readings Page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Correct two existing date values (year typo fixes: +365/+366 days) ---
$ws.Range("B647").Value = 45690.5
$ws.Range("B678").Value = 45721.690972222219

# --- 2) Append new readings as rows 713-720 ---

# Row 713
$ws.Range("A713").Value = "2024-25"
$ws.Range("B713").Value = 45756.449305555558
$ws.Range("C713").Value = 1.94

# Row 714
$ws.Range("A714").Value = "2024-25"
$ws.Range("B714").Value = 45757.505555555559
$ws.Range("C714").Value = 1.94
$ws.Range("D714").Value = 0.03
$ws.Range("E714").Value = "Rain (not SWE)"

# Row 715
$ws.Range("A715").Value = "2024-25"
$ws.Range("B715").Value = 45758.6875
$ws.Range("C715").Value = 1.94
$ws.Range("D715").Value = 0.01
$ws.Range("E715").Value = "Rain (not SWE)"

# Row 716
$ws.Range("A716").Value = "2024-25"
$ws.Range("B716").Value = 45759.540277777778
$ws.Range("C716").Value = 1.94
$ws.Range("E716").Value = "Flynn lake ice out"

# Row 717
$ws.Range("A717").Value = "2024-25"
$ws.Range("B717").Value = 45760.499305555553
$ws.Range("C717").Value = 1.94
$ws.Range("D717").Value = 0.08
$ws.Range("E717").Value = "Rain (not SWE)"

# Row 718
$ws.Range("A718").Value = "2024-25"
$ws.Range("B718").Value = 45761.541666666664
$ws.Range("C718").Value = 1.96
$ws.Range("D718").Value = 0.09
$ws.Range("E718").Value = "Rain (not SWE)"

# Row 719
$ws.Range("A719").Value = "2025-26"
$ws.Range("B719").Value = 45762.522222222222
$ws.Range("C719").Value = 2
$ws.Range("D719").Value = 0.37
$ws.Range("E719").Value = "Mixed (so SWE); Official Eagle lake ice out"

# Row 720
$ws.Range("A720").Value = "2025-26"
$ws.Range("B720").Value = 45763.513194444444
$ws.Range("C720").Value = 1.98

# --- 3) Apply date/number formatting consistent with the rest of column B/C ---
$ws.Range("B713:B720").NumberFormat = $ws.Range("B712").NumberFormat
$ws.Range("C713:C720").NumberFormat = $ws.Range("C712").NumberFormat

# --- 4) Update frozen-pane / selection view state to reflect new bottom rows ---
$ws.Application.ActiveWindow.ScrollRow = 693
$ws.Range("A720").Select()
